# Updates "想去人数" (want-to-go count, column F) figures across the four
# sheets of the 广州-漫展信息 workbook, plus a title tweak on the
# "KANAKO ITO&AYANE" live show rows (drops the "【大会员抢先购】" prefix
# now that the early-access window has closed).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ("Exhibitions") ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value  = 1504
$ws1.Range("F4").Value  = 787
$ws1.Range("F5").Value  = 197
$ws1.Range("F6").Value  = 51
$ws1.Range("F7").Value  = 1091
$ws1.Range("F8").Value  = 668
$ws1.Range("F9").Value  = 752
$ws1.Range("F10").Value = 1319
$ws1.Range("F11").Value = 267
$ws1.Range("F12").Value = 1004
$ws1.Range("F13").Value = 14
$ws1.Range("F14").Value = 52
$ws1.Range("F16").Value = 37
$ws1.Range("F17").Value = 404
$ws1.Range("F19").Value = 288
$ws1.Range("F20").Value = 526
$ws1.Range("F21").Value = 544
$ws1.Range("F23").Value = 218
$ws1.Range("F24").Value = 158

# --- Sheet 2: 演出 ("Performances") ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F3").Value  = 975
$ws2.Range("F5").Value  = 223
$ws2.Range("F8").Value  = 58
$ws2.Range("F9").Value  = 574
$ws2.Range("C10").Value = "广州·KANAKO ITO&AYANE 2024 LIVE"
$ws2.Range("F10").Value = 38
$ws2.Range("F11").Value = 11

# --- Sheet 3: 本地生活 ("Local life") ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 171

# --- Sheet 4: 全部类型 ("All categories") ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value  = 171
$ws4.Range("F4").Value  = 1504
$ws4.Range("F6").Value  = 787
$ws4.Range("F7").Value  = 197
$ws4.Range("F8").Value  = 975
$ws4.Range("F9").Value  = 51
$ws4.Range("F10").Value = 1091
$ws4.Range("F11").Value = 668
$ws4.Range("F12").Value = 752
$ws4.Range("F13").Value = 1319
$ws4.Range("F14").Value = 267
$ws4.Range("F15").Value = 1004
$ws4.Range("F16").Value = 14
$ws4.Range("F17").Value = 52
$ws4.Range("F19").Value = 37
$ws4.Range("F20").Value = 404
$ws4.Range("F22").Value = 223
$ws4.Range("F24").Value = 288
$ws4.Range("F28").Value = 526
$ws4.Range("F29").Value = 544
$ws4.Range("F31").Value = 218
$ws4.Range("F32").Value = 58
$ws4.Range("F33").Value = 158
$ws4.Range("F34").Value = 574
$ws4.Range("C35").Value = "广州·KANAKO ITO&AYANE 2024 LIVE"
$ws4.Range("F35").Value = 38
$ws4.Range("C36").Value = "广州·KANAKO ITO&AYANE 2024 LIVE"
$ws4.Range("F36").Value = 38
$ws4.Range("F37").Value = 11
